# "Charges" (B5) drops from 5% to 0.5%, which ripples through the
# Return / Compounding rate / projected balance formulas below it and
# widens the chart's y-range, hence the x-axis layout fix in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B5").Value = 0.005

# Leave the selection where the author left it when they saved.
$ws.Range("C6").Select()
